$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-11 07:58:05"
$wsZh.Range("G2").Value = "2016-01-11 07:59:11"

# de-de sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-11 07:58:22"
$wsDe.Range("G2").Value = "2016-01-11 07:59:40"
